$wb = $excel.ActiveWorkbook
try {
  $wb.OpenLinks("/valuationquan/HwabaoWPszseinnovation100ETF.xlsx")
} catch {
  Write-Host "ERR: $_"
}
Write-Host "done"
